## "updates from other projects"
##
## 1) Refresh the cached "datetimeFigureOut" footer-date text on every
##    slide layout's Date Placeholder (2/8/2017 -> 1/25/2019).
## 2) Recolor seven "Straight Arrow Connector" line shapes on slide 1
##    from theme accent3 (lumMod 75%) to the fixed RGB 3E4D1F.
## 3) Re-stack "Straight Arrow Connector 43" (id 44) to the very top of
##    the z-order (last shape in the tree) on slide 1.

$p = $ppt.ActivePresentation

# --- 1) Date placeholder text on every slide layout -----------------
$master = $p.SlideMaster
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "1/25/2019"
        }
    }
}

# --- 2) / 3) Connector line-color + z-order fix on slide 1 ----------
$s = $p.Slides.Item(1)

# ids of the connectors whose line color changes from
# schemeClr accent3 (lumMod 75000) to srgbClr 3E4D1F
$recolorIds = @(15, 23, 36, 115, 238, 239, 52)

$moveId = 44
$moveShape = $null

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($recolorIds -contains $sh.Id) {
        $sh.Line.ForeColor.RGB = 2051390   # RGB(0x1F, 0x4D, 0x3E) == &H3E4D1F
    }
    if ($sh.Id -eq $moveId) {
        $moveShape = $sh
    }
}

if ($moveShape -ne $null) {
    $moveShape.ZOrder(0)   # msoBringToFront - send to the end of the shape tree
}
